# Update DateBase/orders/International Ever Green_2025-10-15.xlsx
# Orders sheet: rows 2-8 get new flower data, rows 9-10 are newly added.
# Summary sheet: G2 tracking/order code changes.

$wb = $excel.ActiveWorkbook
$orders = $wb.Worksheets.Item("Orders")
$summary = $wb.Worksheets.Item("Summary")

# The source data keeps every cell (even numeric-looking ones like "36" or
# "0510101055401025") as TEXT, not numbers, so leading zeros / the original
# "number-stored-as-text" shape survive. A leading apostrophe is the normal
# Excel way to force a numeric-looking literal to stay text, so use that for
# the purely-numeric entries instead of letting Excel auto-convert them.

# --- Orders sheet: PackageID (A), FlowerName (C), Number (F) ---

$orders.Range("A2").Value = "'36"

$orders.Range("C2").Value = "574_迷你菊白_undefined_undefined_1bunch"
$orders.Range("F2").Value = "'5"

$orders.Range("C3").Value = "576_迷你菊紫_undefined_undefined_1bunch"
$orders.Range("F3").Value = "'10"

$orders.Range("C4").Value = "575_迷你菊深粉_undefined_undefined_1bunch"
$orders.Range("F4").Value = "'10"

$orders.Range("C5").Value = "384_奶油向日葵_sunflower cream_undefined_1bunch"
$orders.Range("F5").Value = "'10"

$orders.Range("C6").Value = "478_绿芯向日葵_sunflower mini_undefined_1bunch"
$orders.Range("F6").Value = "'5"

$orders.Range("C7").Value = "411_紫罗兰白_violet white_undefined_1bunch"
$orders.Range("F7").Value = "'5"

$orders.Range("C8").Value = "505_紫罗兰紫_violet purple_undefined_1bunch"
$orders.Range("F8").Value = "'40"

# New rows 9 and 10
$orders.Range("C9").Value = "506_紫罗兰香槟色_violet champagne_undefined_1bunch"
$orders.Range("F9").Value = "'10"

$orders.Range("C10").Value = "412_紫罗兰粉_violet pink_undefined_1bunch"
$orders.Range("F10").Value = "'25"

# --- Summary sheet: tracking/order code ---
$summary.Range("G2").Value = "'0510101055401025"
